$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-89)
# from 45185 (2023-09-16) to 45204 (2023-10-05).
$ws.Range("C2:C89").Value = 45204
